$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1924.5
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H122").Value = 1924.5
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()
$ws.Range("H131").Value = 13462.125
$ws.Range("I131").Value = 13462.125
$ws.Range("K131").Value = 40386.375
$ws.Range("M131").Value = -35346.375
$ws.Range("H141").Value = 2832
$ws.Range("I141").Value = 2832
$ws.Range("K141").Value = 8496
$ws.Range("M141").Value = -3316
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 138
$ws.Range("J5").Value = 500
$ws.Range("L5").Value = 500
$ws.Range("N5").Value = -724
$ws.Range("H32").Value = 3211493.5
$ws.Range("I32").Value = 3336849.8
$ws.Range("J32").Value = 2333999.8
$ws.Range("K32").Value = 3336849.8
$ws.Range("L32").Value = 2333999.8
$ws.Range("M32").Value = -3336562.8
$ws.Range("N32").Value = -2334573.8
$ws.Range("H62").Value = 39000
$ws.Range("J62").Value = 39000
$ws.Range("L62").Value = 39000
$ws.Range("N62").Value = -40248
$ws.Range("H65").Value = 39000
$ws.Range("J65").Value = 39000
$ws.Range("L65").Value = 117000
$ws.Range("N65").Value = -123240
$ws.Range("H97").Value = 867.06665
$ws.Range("J97").Value = 1442.5
$ws.Range("L97").Value = 1442.5
$ws.Range("N97").Value = -2434.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 138
$ws.Range("J4").Value = 500
$ws.Range("L4").Value = 500
$ws.Range("N4").Value = -730
$ws.Range("H99").Value = 2148.8333
$ws.Range("I99").Value = 2487.5
$ws.Range("J99").Value = 1979.5
$ws.Range("K99").Value = 2487.5
$ws.Range("L99").Value = 1979.5
$ws.Range("M99").Value = -989.5
$ws.Range("N99").Value = -4975.5
$ws.Range("H105").Value = 1899.2
$ws.Range("I105").Value = 1925
$ws.Range("K105").Value = 1925
$ws.Range("M105").Value = -178
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 830.625
$ws.Range("I22").Value = 724.25
$ws.Range("K22").Value = 724.25
$ws.Range("M22").Value = -374.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 110.833336
$ws.Range("I17").Value = 110
$ws.Range("J17").Value = 112.5
$ws.Range("K17").Value = 330
$ws.Range("L17").Value = 337.5
$ws.Range("M17").Value = -161
$ws.Range("N17").Value = -675.5
$ws.Range("H22").Value = 2165.5789
$ws.Range("I22").Value = 2000
$ws.Range("J22").Value = 2174.7778
$ws.Range("K22").Value = 6000
$ws.Range("L22").Value = 6524.3334
$ws.Range("M22").Value = -5831
$ws.Range("N22").Value = -6862.3334
$ws.Range("H27").Value = 2165.5789
$ws.Range("I27").Value = 2000
$ws.Range("J27").Value = 2174.7778
$ws.Range("K27").Value = 6000
$ws.Range("L27").Value = 6524.3334
$ws.Range("M27").Value = -5898
$ws.Range("N27").Value = -6728.3334
$ws.Range("H137").Value = 3652.75
$ws.Range("I137").Value = 2395
$ws.Range("J137").Value = 4072
$ws.Range("K137").Value = 7185
$ws.Range("L137").Value = 12216
$ws.Range("M137").Value = -2085
$ws.Range("N137").Value = -22416
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3599.75
$ws.Range("J80").Value = 3600
$ws.Range("L80").Value = 3600
$ws.Range("N80").Value = -5596
$ws.Range("H83").Value = 3599.75
$ws.Range("J83").Value = 3600
$ws.Range("L83").Value = 18000
$ws.Range("N83").Value = -27984
$ws.Range("H107").Value = 1275
$ws.Range("I107").Value = 318.85715
$ws.Range("J107").Value = 2948.25
$ws.Range("K107").Value = 318.85715
$ws.Range("L107").Value = 2948.25
$ws.Range("M107").Value = 1601.14285
$ws.Range("N107").Value = -6788.25
$ws.Range("H126").Value = 8500
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 5861
$ws.Range("I132").Value = 6092.9287
$ws.Range("J132").Value = 2614
$ws.Range("K132").Value = 18278.7861
$ws.Range("L132").Value = 7842
$ws.Range("M132").Value = -15748.7861
$ws.Range("N132").Value = -12902
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1699.8
$ws.Range("J16").Value = 500
$ws.Range("L16").Value = 500
$ws.Range("N16").Value = -840
$ws.Range("H22").Value = 2253.9167
$ws.Range("I22").Value = 2487.5
$ws.Range("J22").Value = 1786.75
$ws.Range("K22").Value = 2487.5
$ws.Range("L22").Value = 1786.75
$ws.Range("M22").Value = -2192.5
$ws.Range("N22").Value = -2376.75
$ws.Range("H27").Value = 2253.9167
$ws.Range("I27").Value = 2487.5
$ws.Range("J27").Value = 1786.75
$ws.Range("K27").Value = 2487.5
$ws.Range("L27").Value = 1786.75
$ws.Range("M27").Value = -2380.5
$ws.Range("N27").Value = -2000.75
$ws.Range("H100").Value = 4256.3335
$ws.Range("J100").Value = 2798.6667
$ws.Range("L100").Value = 2798.6667
$ws.Range("N100").Value = -3880.6667
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 34666.332
$ws.Range("I62").Value = 28999
$ws.Range("J62").Value = 37500
$ws.Range("K62").Value = 28999
$ws.Range("L62").Value = 37500
$ws.Range("M62").Value = -28375
$ws.Range("N62").Value = -38748
$ws.Range("H65").Value = 34666.332
$ws.Range("I65").Value = 28999
$ws.Range("J65").Value = 37500
$ws.Range("K65").Value = 144995
$ws.Range("L65").Value = 187500
$ws.Range("M65").Value = -141875
$ws.Range("N65").Value = -193740
$ws.Range("H100").Value = 6667319.5
$ws.Range("I100").Value = 7692873
$ws.Range("J100").Value = 1222
$ws.Range("K100").Value = 15385746
$ws.Range("L100").Value = 2444
$ws.Range("M100").Value = -15385205
$ws.Range("N100").Value = -3526
$ws.Range("H132").Value = 1124.9231
$ws.Range("I132").Value = 1124.9231
$ws.Range("K132").Value = 3374.7693
$ws.Range("M132").Value = -844.7692999999999
